$wb = $excel.ActiveWorkbook

# --- 1. Fix the date typo in the shared string "11/11/2020" -> "11/11/2021" ---
# Only X3:X25 get the style change further below, but the text fix applies to
# every cell sharing that string (X2:X25).
$lots = $wb.Worksheets.Item("lots")
$used = $lots.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $lots.Cells.Item($r, 24)
    if ($cell.Value -eq "11/11/2020") {
        $cell.Value = "11/11/2021"
    }
}

# --- 2. Give X3:X25 a distinct (but visually identical) font/style ---
$range = $lots.Range("X3:X25")
$range.Font.Name = "Calibri"
$range.Font.Size = 11
$range.Font.Color = 0
$range.Font.Charset = 0

# --- 3. Update the remembered selection on each sheet ---
$lots.Select()
$lots.Range("X26").Select()

foreach ($name in @("MatieresPremieres", "Biocarburants", "Pays", "Societes", "SitesDeLivraison")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Select()
    $sheet.Range("A1").Select()
}

$lots.Select()
